# Add a new "2019-08-26" sheet in front of the existing sheets, seeded
# from the current "2019-08-19" sheet's data (same neighborhoods / same
# "ongoing transmission" and "all cases" counts), refresh the header
# labels for the new reporting date, and apply the two "latest cases"
# corrections (Williamsburg and TOTAL drop from 1 to 0).
#
# NOTE: object references returned by Worksheets.Item(...) are
# index-anchored, so they go stale as soon as the sheet collection is
# reshuffled (e.g. by Worksheets.Add). We therefore re-resolve every
# worksheet reference by name immediately before using it, rather than
# holding on to a reference captured before the Add() call.

$wb = $excel.ActiveWorkbook

# Insert the new sheet immediately before "2019-08-19" (i.e. as the
# first tab), then rename it and re-fetch stable references by name.
$sourceForPlacement = $wb.Worksheets.Item("2019-08-19")
$placeholder = $wb.Worksheets.Add($sourceForPlacement)
$placeholder.Name = "2019-08-26"

$newSheet = $wb.Worksheets.Item("2019-08-26")
$source = $wb.Worksheets.Item("2019-08-19")

# Copy the neighborhood names and the "ongoing transmission" / "all
# cases" columns verbatim from the 2019-08-19 sheet; the "latest cases"
# column (D) is recomputed below.
for ($r = 1; $r -le 20; $r++) {
    for ($c = 1; $c -le 4; $c++) {
        $val = $source.Cells.Item($r, $c).Value2
        $newSheet.Cells.Item($r, $c).Value = $val
    }
}

# Refresh the header row for the new reporting date.
$newSheet.Cells.Item(1, 1).Value = "neighborhood"
$newSheet.Cells.Item(1, 2).Value = "ongoing transmission (2019-08-26)"
$newSheet.Cells.Item(1, 3).Value = "all cases (2018-09-01 to 2019-08-26)"
$newSheet.Cells.Item(1, 4).Value = "latest cases (2019-08-19 to 2019-08-26)"

# Apply the two "latest cases" corrections noted for this edition:
# Williamsburg (row 4) and TOTAL (row 20) both drop from 1 to 0.
$newSheet.Cells.Item(4, 4).Value = 0
$newSheet.Cells.Item(20, 4).Value = 0

# Match the column widths used on the other sheets.
$newSheet.Columns.Item(1).ColumnWidth = 26.32
$newSheet.Columns.Item(2).ColumnWidth = 29.16
$newSheet.Columns.Item(3).ColumnWidth = 30.86
$newSheet.Columns.Item(4).ColumnWidth = 33.38

# Put the selection back at A1 on both the new sheet and the sheet that
# used to be first, then make the new sheet the active tab.
$newSheet.Range("A1").Select() | Out-Null
$source.Range("A1").Select() | Out-Null
$newSheet.Activate() | Out-Null
$newSheet.Range("A1").Select() | Out-Null
